$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Worksheet, $CellRef, $TextValue)
    $cell = $Worksheet.Range($CellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $TextValue
    $cell.Style = "Normal"
}

Set-TextValue $ws "D2" "26.174.34"
Set-TextValue $ws "E2" "  -2.16%  "
Set-TextValue $ws "D3" "1.670.66"
Set-TextValue $ws "E3" "  -1.63%  "
Set-TextValue $ws "E4" "  +0.13%  "
Set-TextValue $ws "D5" "217.23"
Set-TextValue $ws "E5" "  -0.96%  "
Set-TextValue $ws "D6" "0.5115"
Set-TextValue $ws "E6" "  +0.39%  "
Set-TextValue $ws "E7" "  +0.15%  "
Set-TextValue $ws "D8" "0.2640"
Set-TextValue $ws "E8" "  +1.28%  "
Set-TextValue $ws "D9" "0.06394"
Set-TextValue $ws "E9" "  +3.94%  "
Set-TextValue $ws "D10" "21.59"
Set-TextValue $ws "E10" "  -1.65%  "
Set-TextValue $ws "D11" "0.07414"
Set-TextValue $ws "E11" "  +0.97%  "
Set-TextValue $ws "D12" "1.677.03"
Set-TextValue $ws "E12" "  -0.60%  "
Set-TextValue $ws "D13" "4.526"
Set-TextValue $ws "D14" "0.5812"
Set-TextValue $ws "E14" "  +1.10%  "
Set-TextValue $ws "D15" "0.000008563"
Set-TextValue $ws "E15" "  +4.30%  "
Set-TextValue $ws "D16" "64.38"
Set-TextValue $ws "E16" "  -2.02%  "
Set-TextValue $ws "D17" "26.241.20"
Set-TextValue $ws "E17" "  -2.04%  "
Set-TextValue $ws "D18" "4.931"
Set-TextValue $ws "E18" "  -1.91%  "
Set-TextValue $ws "E19" "  +0.07%  "
Set-TextValue $ws "D20" "10.84"
Set-TextValue $ws "E20" "  +0.97%  "
Set-TextValue $ws "D21" "189.98"
Set-TextValue $ws "E21" "  +2.41%  "
Set-TextValue $ws "E22" "  -0.46%  "
Set-TextValue $ws "D23" "1.007"
Set-TextValue $ws "E23" "  +0.09%  "
Set-TextValue $ws "D24" "145.39"
Set-TextValue $ws "E24" "  -0.15%  "
Set-TextValue $ws "D25" "7.627"
Set-TextValue $ws "E25" "  -0.75%  "
Set-TextValue $ws "D26" "0.1186"
Set-TextValue $ws "E26" "  +3.21%  "
Set-TextValue $ws "D27" "15.66"
Set-TextValue $ws "E27" "  +2.24%  "
Set-TextValue $ws "D28" "0.06381"
Set-TextValue $ws "E28" "  +12.18%  "
Set-TextValue $ws "D29" "1.300"
Set-TextValue $ws "E29" "  -1.31%  "
Set-TextValue $ws "D30" "1.320"
Set-TextValue $ws "E30" "  -1.41%  "
Set-TextValue $ws "D31" "3.534"
Set-TextValue $ws "E31" "  +1.38%  "
Set-TextValue $ws "D32" "3.525"
Set-TextValue $ws "E32" "  +1.85%  "
Set-TextValue $ws "D33" "1.639"
Set-TextValue $ws "E33" "  -1.60%  "
Set-TextValue $ws "E34" "  +0.92%  "
Set-TextValue $ws "D35" "0.6084"
Set-TextValue $ws "E35" "  +2.70%  "
Set-TextValue $ws "D36" "2.375"
Set-TextValue $ws "E36" "  -1.50%  "
Set-TextValue $ws "D37" "2.657"
Set-TextValue $ws "E37" "  +0.55%  "
Set-TextValue $ws "D38" "6.159"
Set-TextValue $ws "E38" "  +3.36%  "
Set-TextValue $ws "D39" "0.01604"
Set-TextValue $ws "E39" "  +0.66%  "
Set-TextValue $ws "D40" "1.082.96"
Set-TextValue $ws "E40" "  +1.16%  "
Set-TextValue $ws "D41" "0.8655"
Set-TextValue $ws "E42" "  +0.72%  "
Set-TextValue $ws "D43" "101.23"
Set-TextValue $ws "E43" "  +2.56%  "
Set-TextValue $ws "D44" "1.819.65"
Set-TextValue $ws "E44" "  -1.82%  "
Set-TextValue $ws "D45" "0.00000000112"
Set-TextValue $ws "E45" "  +5.85%  "
Set-TextValue $ws "D46" "56.25"
Set-TextValue $ws "E46" "  -0.76%  "
Set-TextValue $ws "D47" "1.008"
Set-TextValue $ws "E47" "  -0.03%  "
Set-TextValue $ws "D48" "8.099"
Set-TextValue $ws "E48" "  +0.99%  "
Set-TextValue $ws "E49" "  -0.15%  "
Set-TextValue $ws "D50" "0.4297"
Set-TextValue $ws "E50" "  -0.85%  "
Set-TextValue $ws "B51" "RenderToken"
Set-TextValue $ws "C51" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws "D51" "1.441"
Set-TextValue $ws "E51" "  -0.70%  "
